# Applies the 8800008.xlsx content update:
#  - "Objetivos:" answer cell replaced with the professor name string
#  - "Docentes responsáveis:" answer row removed (row 13 repurposed)
#  - "Programa resumido:" answer replaced with "Semestral"
#  - "Programa:" long-text answer removed; the old "Ativação:" date value
#    ends up echoed under "Programa:" (matches source data exactly)
#  - "Método:" row now carries the professor-name string instead of being
#    the criteria-table header; "Critério:", "Norma de recuperação:" and
#    "Bibliografia:" labels shift up one row, and the old Bibliografia
#    long-text row (old row 22) is deleted entirely
#  - row heights adjusted to match the new row contents

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last row (old "Bibliografia:" long-text row); this shrinks the
# sheet's dimension from C22 to C21 and leaves rows 1-21 untouched.
$ws.Rows.Item(22).Delete()

# Row 10 ("Objetivos:"): answer becomes the professor-name string.
# B10/C10 already exist, so they keep their normal/red-text column styles.
$ws.Range("B10").Value = "198273 - Domingos Savio Giordani"
$ws.Range("C10").Value = "198273 - Domingos Savio Giordani"

# Row 13: now holds "Programa resumido:" / "Semestral" (was the bare
# "Docentes responsáveis:" answer row, A13 previously empty).
# B13/C13 already exist, so they keep their normal/red-text column styles.
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Rows.Item(13).RowHeight = 60

# Row 14: label becomes "Short syllabus:"; its old answer text is removed.
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14:C14").Clear()

# Row 15: label becomes "Programa:"; answer becomes the activation-date
# text, and the row grows to the 120pt "long answer" height. B15/C15
# don't exist yet on this row, so two things need help:
#  1) writing a date-shaped string directly would get auto-parsed into a
#     date serial number instead of staying literal text;
#  2) a freshly-created cell in column B would inherit the wrong style
#     (the sheet's <cols> ranges for columns A and B overlap, so a brand
#     new cell resolves to the column-A style instead of column B's).
# Paste (values-only, then formats-only) from B8/C8 - which already hold
# "01/01/2020" as literal text with the correct column style - to
# materialize B15/C15 correctly without touching styles.xml at all.
$ws.Range("A15").Value = "Programa:"
$ws.Range("B8").Copy()
$ws.Range("B15").PasteSpecial(-4163)
$ws.Range("B8").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("C8").Copy()
$ws.Range("C15").PasteSpecial(-4163)
$ws.Range("C8").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Rows.Item(15).RowHeight = 120

# Row 16: label becomes "Syllabus:"; its old long-text answer is removed.
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16:C16").Clear()

# Row 17: label becomes "Avaliação:"; reverts to the default row height.
$ws.Range("A17").Value = "Avaliação:"
$ws.Rows.Item(17).EntireRow.AutoFit()

# Row 18: label becomes "Método:"; answer becomes the professor-name
# string, and the row grows to the 60pt "answer" height. B18/C18 don't
# exist yet either, so materialize them (values, then formats) from
# B19/C19 - same normal/red-text column style - before writing the final
# text, avoiding any new styles.xml entries.
$ws.Range("A18").Value = "Método:"
$ws.Range("B19").Copy()
$ws.Range("B18").PasteSpecial(-4163)
$ws.Range("B19").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("B18").Value = "198273 - Domingos Savio Giordani"
$ws.Range("C19").Copy()
$ws.Range("C18").PasteSpecial(-4163)
$ws.Range("C19").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("C18").Value = "198273 - Domingos Savio Giordani"
$ws.Rows.Item(18).RowHeight = 60

# Row 19: label becomes "Critério:" (answer text unchanged).
$ws.Range("A19").Value = "Critério:"

# Row 20: label becomes "Norma de recuperação:" (answer text unchanged).
$ws.Range("A20").Value = "Norma de recuperação:"

# Row 21: label becomes "Bibliografia:" (answer text unchanged); grows to
# the 120pt "long answer" height.
$ws.Range("A21").Value = "Bibliografia:"
$ws.Rows.Item(21).RowHeight = 120

Write-Host "applied 8800008.xlsx content update"
